# Update "98_1" confirmations summary sheet:
#  - Remove the standalone "Summary" section header row (old row 28), shifting
#    the totals rows up by one.
#  - Prefix each branch's sub-category labels with the branch/category name
#    (e.g. "     New nominations" -> "     Army, New nominations").
#  - Rename "Civilian nominations" header to "Civilian".
#  - Fix wording/typos in the total labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Summary" header row; everything below shifts up one row.
$ws.Rows(28).Delete()

# Army section (row 6 header, rows 7-10 detail)
$ws.Range("A7").Value  = "     Army, New nominations"
$ws.Range("A8").Value  = "     Army, Confirmed "
$ws.Range("A9").Value  = "     Army, Failed at August-September adjournment"
$ws.Range("A10").Value = "     Army, Failed at November sine die adjournment"

# Navy section (row 11 header, rows 12-13 detail)
$ws.Range("A12").Value = "     Navy, New nominations"
$ws.Range("A13").Value = "     Navy, Confirmed "

# Air Force section (row 14 header, rows 15-18 detail)
$ws.Range("A15").Value = "     Air Force, New nominations"
$ws.Range("A16").Value = "     Air Force, Confirmed "
$ws.Range("A17").Value = "     Air Force, Unconfirmed "
$ws.Range("A18").Value = "     Air Force, Failed at August-September adjournment"

# Marine Corps section (row 19 header, rows 20-21 detail)
$ws.Range("A20").Value = "     Marine Corps, New nominations"
$ws.Range("A21").Value = "     Marine Corps, Confirmed "

# Civilian section (row 22 header, rows 23-27 detail)
$ws.Range("A22").Value = "Civilian"
$ws.Range("A23").Value = "     Civilian, New nominations"
$ws.Range("A24").Value = "     Civilian, Confirmed "
$ws.Range("A25").Value = "     Civilian, Withdrawn "
$ws.Range("A26").Value = "     Civilian, Failed at August-September adjournment"
$ws.Range("A27").Value = "     Civilian, Failed at November sine die adjournment"

# Totals section (now rows 28-33 after the Summary row was removed)
$ws.Range("A28").Value = "Total new nominations"
$ws.Range("A29").Value = "Total confirmed "
$ws.Range("A33").Value = "Total failed at November sine die adjournment "
